$d = $word.ActiveDocument

# Locate the "Video hint" list paragraph that immediately precedes the
# "TICKABLE It is very easy to get help in Sage..." paragraph. (The text
# "Video hint" recurs several times in the document as the caption of a
# hyperlink before each TICKABLE exercise, so we disambiguate using the
# paragraph that follows it.)
$target = $null
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count - 1; $i++) {
    $p = $paragraphs.Item($i)
    $pText = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($pText -eq "Video hint") {
        $nextText = $paragraphs.Item($i + 1).Range.Text
        if ($nextText.StartsWith("TICKABLE It is very easy to get help in Sage")) {
            $target = $p
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Video hint' paragraph preceding the help TICKABLE item."
}

# Insert a brand new paragraph right after it, inheriting the same bullet
# list (numId) as "Video hint" so it slots into the list naturally.
$target.Range.InsertParagraphAfter()

$newPara = $paragraphs.Item($i + 1)
$newRange = $newPara.Range

$message = "If you forget your password DO NOT CREATE ANOTHER ACCOUNT: come and see me (Vince Knight) and I can reset your password."
$newRange.InsertAfter($message)

# Bold only the literal text, not the trailing paragraph mark.
$textRange = $d.Range($newRange.Start, $newRange.Start + $message.Length)
$textRange.Bold = 1

Write-Output "Inserted password-reset paragraph after 'Video hint' (paragraph index $i)."
